# daily auto push: 2026-02-03 14:09 UTC
# Two new hourly readings for 2026/02/03 were appended to the data set
# (continuing the 7 / 13 / 16 sequence already present in rows 772-774),
# pushing all the later rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new data points by inserting two blank rows
# right above the existing "2026/12/29" block (old row 775).
$ws.Range("A775:D776").Insert()

# New row 775: 2026/02/03, 火, hour=19, value=201
$ws.Cells.Item(775, 1).NumberFormat = "@"
$ws.Cells.Item(775, 1).Value = "2026/02/03"
$ws.Cells.Item(775, 1).ClearFormats()
$ws.Cells.Item(775, 2).Value = "火"
$ws.Cells.Item(775, 3).Value = 19
$ws.Cells.Item(775, 4).Value = 201

# New row 776: 2026/02/03, 火, hour=22, value=201
$ws.Cells.Item(776, 1).NumberFormat = "@"
$ws.Cells.Item(776, 1).Value = "2026/02/03"
$ws.Cells.Item(776, 1).ClearFormats()
$ws.Cells.Item(776, 2).Value = "火"
$ws.Cells.Item(776, 3).Value = 22
$ws.Cells.Item(776, 4).Value = 201
